$wb = $excel.ActiveWorkbook

# --- Skills_GET: user moved the selection to D5 ---
$wsGet = $wb.Worksheets.Item("Skills_GET")
$wsGet.Activate()
$wsGet.Range("D5").Select()

# --- Skills_POST: a few input values changed, plus the status message/style
#     on row 8 was swapped to match the "success" row (row 2) ---
$wsPost = $wb.Worksheets.Item("Skills_POST")
$wsPost.Activate()
$wsPost.Range("B2").Value = 36
$wsPost.Range("B8").Value = 34
$wsPost.Range("D8").Value = 201
$wsPost.Range("E8").Value = "Skill record successfully created"

# Copy E2's formatting (style index) onto E8 so it matches the "success" row
$wsPost.Range("E2").Copy()
$wsPost.Range("E8").PasteSpecial(-4122)   # xlPasteFormats

$wsPost.Range("B9").Value = 35
$wsPost.Range("B9").Select()
$excel.ActiveWindow.ScrollRow = 2
$excel.ActiveWindow.ScrollColumn = 2

# --- Skills_PUT: selection moved to D8, view scrolled so row 4 is on top ---
$wsPut = $wb.Worksheets.Item("Skills_PUT")
$wsPut.Activate()
$wsPut.Range("D8").Select()
$excel.ActiveWindow.ScrollRow = 4
$excel.ActiveWindow.ScrollColumn = 1

# --- Skills_DELETE: selection moved to B2, and the Skill_Id input changed ---
$wsDelete = $wb.Worksheets.Item("Skills_DELETE")
$wsDelete.Activate()
$wsDelete.Range("B2").Value = 33
$wsDelete.Range("B2").Select()

# Leave Skills_PUT as the active sheet/tab (it was tabSelected in the file)
$wsPut.Activate()
